$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13. Everything that was row 13 ("Programa
# resumido:" ...) and below shifts down to row 14+, carrying its row
# heights along automatically.
$ws.Rows(13).Insert()

# --- Row 10 (Objetivos:) ---------------------------------------------------
# B/C previously (incorrectly) duplicated the teacher name; replace with the
# real objectives paragraph.
$objetivos = "Levar aos estudantes conhecimentos básicos sobre: a) processos fermentativos, com ênfase em processos de interesse industrial; b) bioquímica das fermentações focando as rotas metabólicas utilizadas por microrganismos de interesse industrial; c) suas respectivas aplicações em processos industriais, permitindo a determinação de parâmetros de avaliação de desempenho."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# --- Row 13 (new row: Docentes responsáveis: data) -------------------------
$ws.Range("A13").Clear()
$nome = "3403572 - Ismael Maciel de Mancilha"
$ws.Range("B13").Value = $nome
$ws.Range("C13").Value = $nome
# Match the normal / red wrap-text formatting used by the rest of the sheet
# (row insert defaults to the bold "label" style picked up from column A).
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 14 (Programa resumido:) --------------------------------------------
$resumido = "Biotecnologia; processos fermentativos; bioquímica das fermentações (vias metabólicas de interesse industrial); processos fermentativos de interesse industrial."
$ws.Range("B14").Value = $resumido
$ws.Range("C14").Value = $resumido

# --- Row 16 (Programa:) ------------------------------------------------------
$programa = "1. Biotecnologia: conceitos, áreas de aplicação, caráter multidisciplinar e exemplos de produtos biotecnológicos.2. Processos Fermentativos: conceito, exemplos, fases de um processo fermentativo. Modalidades de Processos Fermentativos: a)formas de condução; b) fermentação induzida e espontânea; c) estado físico do meio de fermentação; d) suprimento de oxigênio; e) processos submersos e em superfície; f) cinética de formação de produto em relação do metabolismo primário.3. Bioquímica das fermentações: fermentação  conceitos, objetivos, aerobiose x anaerobiose; balanço energético; estágios preliminares da fermentação (hidrólise extracelular e permeabilidade da membrana); vias metabólica de interesse industrial: a) via glicolítica: reações e controle; fermentação alcoólica, homoláctica, acetona/butanol, ácido-mista e 2,3 butanodiol; b) Via Fosfo-Cetolase: fermentação heteroláctica e c) via Entner Doudoroff: fermentação alcoólica por Zymomonas mobilis. Balanço da Fermentação: % de carbono recuperado e balanço de oxi-redução; parâmetros de avaliação - rendimento, eficiência e produtividade de processos fermentativos; Processos de Interesse: processamento de cacau, produção de etanol, alimentos fermentados e outros."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# --- Row 19 (Método:) ---------------------------------------------------------
$metodo = "A avaliação será feita por meio de provas escritas."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# --- Row 20 (Critério:) --------------------------------------------------------
$criterio = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2) / 2"
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# --- Row 21 (Norma de recuperação:) ---------------------------------------------
$norma = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# --- Row 22 (Bibliografia:) -------------------------------------------------------
$biblio = "1. AMERINE, M.A, OUGH,C.S., Methods for analysis of musts and wines. New York: John Wiley & Sons, 1980. `n2. AMORIM, H.V., Fermentação Alcoólica ciência e tecnologia. Piracicaba: Fermentec, 2006.`n3. BORZANI, W., SCHMIDELL, W., LIMA, U.A., AQUARONE, E. Série de Biotecnologia Vol. 1  Fundamentos e Vol. 4 Processos Fermentativos e Enzimáticos. São Paulo: Ed.Edgard Blucher, 2001.`n4. EL-MANSI, E.M.T., BRYCE, C.E.A., DEMAIN, A.L., ALLMAN,A.R. Fermentation Microbiology and Biotechnology. 2ª Ed. New York: CRC Taylor & Francis, 2007."
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio
